# Commodity Lab Item and Department_uuid.xlsx — CM 4.0.3.3 update
#
# Adds two new lookup tables ("InstanceType" and "Stockroom") below the
# existing "Item Name" / "Department Name" tables, and renames the two
# existing table headers to their "(lab)" variants. Cell write order
# mirrors the original authoring sequence (UUID column filled before the
# label column for data rows) so the regenerated shared-string table
# lines up with the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 3: new InstanceType table (rows 35-37) ---
$ws.Cells.Item(35, 1).Value = "InstanceType"
$ws.Cells.Item(35, 2).Value = "UUID"
$ws.Range("A35:B35").Interior.Color = 65535
$ws.Range("A35:B35").Font.Bold = $true

$ws.Cells.Item(36, 2).Value = "fce0b4fc-9402-424a-aacb-f99599e51a9f"
$ws.Cells.Item(36, 1).Value = "Receipt"

$ws.Cells.Item(37, 2).Value = "c264f34b-c795-4576-9928-454d1fa20e09"
$ws.Cells.Item(37, 1).Value = "Distribution"

# --- Section 1: Item Name(lab) header (row 1) ---
# Rename the header text; item data rows (2-11) are unchanged.
$ws.Cells.Item(1, 1).Value = "Item Name(lab)"

# --- Section 2: Department Name(lab) header (row 13) ---
# Rename header text and make it bold like the other section headers
# (reuses the existing bold+fill style already used by row 1).
$ws.Cells.Item(13, 1).Value = "Department Name(lab)"
$ws.Range("A13:B13").Font.Bold = $true

# --- Section 4: new Stockroom table (rows 41-43) ---
$ws.Cells.Item(41, 1).Value = "Stockroom"
$ws.Cells.Item(41, 2).Value = "UUID"
$ws.Range("A41:B41").Interior.Color = 65535
$ws.Range("A41:B41").Font.Bold = $true

$ws.Cells.Item(42, 2).Value = "2741bae2-c5de-43ef-891f-7ec2fd58f442"
$ws.Cells.Item(43, 2).Value = "5452ec3e-2fe1-46de-8a6e-28c6442e4cc0"

$ws.Cells.Item(42, 1).Value = "Lab StockRoom (lab)"
$ws.Cells.Item(43, 1).Value = "Bulk Store (pharmacy)"

# --- Update selection to mirror the saved view state ---
$ws.Range("A46").Select()
